$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

# --- Simple part-number text swaps (Find/Replace, restricted to whole text) ---
$d.Content.Find.Execute("DE92-02588G", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "WB56X34928", 2)

$d.Content.Find.Execute("5304515738", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "WD21X28718", 2)

$d.Content.Find.Execute("WD21X28958", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "W11747577", 2)

$d.Content.Find.Execute("WD05X35098", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "W11095995", 2)

# --- Qty in Cart for the WD05X35098 / W11095995 row: 1 -> 2 (table row 5, col 2) ---
$t.Cell(5, 2).Range.Text = "2"

# --- Collapse rows 6-13 (154853801 .. DC64-00802B) down to nothing, keeping only ---
# --- the final data row (5304524473 / ... ), whose part number becomes W11498796 ---
for ($i = 0; $i -lt 8; $i++) {
    $t.Rows.Item(6).Delete()
}
$t.Cell(6, 1).Range.Text = "W11498796"
